$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove trailing whitespace from the variable-name cells (shared-string
# cleanup described in the commit message: "remove the spacing").
# Clearing the cell first forces Excel to allocate a fresh shared-string
# entry (matching the edit order observed in the original author's
# workbook) rather than mutating the old shared-string slot in place.
$ws.Range("A64").ClearContents()
$ws.Range("A64").Value = "WHQ070"

$ws.Range("A63").ClearContents()
$ws.Range("A63").Value = "WHQ030"

$ws.Range("A19").ClearContents()
$ws.Range("A19").Value = "LBDHDDSI"

$ws.Range("A20").ClearContents()
$ws.Range("A20").Value = "LBDTCSI"

$ws.Range("A21").ClearContents()
$ws.Range("A21").Value = "LBXVIDMS"

$ws.Range("A22").ClearContents()
$ws.Range("A22").Value = "LBXGH"

# Update the view state: scroll position and active cell selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("A59").Select()
